# Fruta / hortaliza, semanal
# Insert a new weekly price-report row before the current row 94 (pushing
# the existing rows 94-133 down to 95-134) and populate it with the new
# week's data for "Macroferia Regional de Talca - Berenjena".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 94:133 down to 95:134, leaving a blank row 94 behind
# (formatting of row 94, e.g. the date style in column D, is inherited
# from the row above, matching native Excel "Insert" behaviour).
$ws.Rows("94:94").Insert()

# Populate the newly inserted row 94 with the new data point.
$ws.Cells.Item(94, 1).Value  = 5
$ws.Cells.Item(94, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(94, 3).Value  = "Maule"
$ws.Cells.Item(94, 4).Value  = 44960
$ws.Cells.Item(94, 5).Value  = 7
$ws.Cells.Item(94, 6).Value  = 100112001
$ws.Cells.Item(94, 7).Value  = "Berenjena"
$ws.Cells.Item(94, 8).Value  = "Sin especificar"
$ws.Cells.Item(94, 9).Value  = "Primera"
$ws.Cells.Item(94, 10).Value = 200
$ws.Cells.Item(94, 11).Value = 6000
$ws.Cells.Item(94, 12).Value = 6000
$ws.Cells.Item(94, 13).Value = 6000
$ws.Cells.Item(94, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(94, 15).Value = "Región del Maule"
$ws.Cells.Item(94, 16).Value = 120
$ws.Cells.Item(94, 17).Value = 50
$ws.Cells.Item(94, 18).Value = "Hortaliza"
